$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log entry as row 40 on the Logs sheet
$logs.Range("A40").Value = "Sponsoraanvraag"
$logs.Range("B40").Value = "mailmind.test@zohomail.eu"
$logs.Range("C40").Value = "Zou uw bedrijf bereid zijn om ons sportevenement te sponsoren?"
$logs.Range("D40").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F40").Value = "2025-06-19 22:22:16"
$logs.Range("G40").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row (D2:D39 -> D2:D40, G2:G39 -> G2:G40)
$fcsD = $logs.Range("D2:D39").FormatConditions
for ($i = 1; $i -le $fcsD.Count; $i++) {
    $fcsD.Item($i).ModifyAppliesToRange($logs.Range("D2:D40"))
}

$fcsG = $logs.Range("G2:G39").FormatConditions
for ($i = 1; $i -le $fcsG.Count; $i++) {
    $fcsG.Item($i).ModifyAppliesToRange($logs.Range("G2:G40"))
}

# Update the Dashboard summary count for "Samenwerking / Partnerverzoek"
$dash.Range("B2").Value = 11
